# Update the "K" column (column G) values on Sheet1 to reflect the
# regenerated save_data (K computed instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 2
    4  = 6
    5  = 2
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 4
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 2
    22 = 3
    23 = 1
    24 = 1
    25 = 1
    26 = 3
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 2
    37 = 0
    38 = 2
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 3
    44 = 1
    45 = 2
    46 = 0
    47 = 1
    48 = 0
    49 = 0
    50 = 1
    51 = 0
    52 = 2
    53 = 3
    54 = 2
    55 = 1
    56 = 3
    57 = 1
    58 = 3
    59 = 0
    60 = 1
    61 = 1
    62 = 2
    64 = 2
    65 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
